# xls export geometry fix, added server power\health status
#
# Relabel a few hwinvent report headers to shorter captions and shrink
# the columns that held the old, longer text so the export reads cleanly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header text relabels (row 1) ---------------------------------------
$ws.Range("F1").Value = "Memory tot.size"     # was "System memory size"
$ws.Range("H1").Value = "Memory P/Ns"         # was "Memory module part number"
$ws.Range("M1").Value = "HDD slot pop."       # was "HDD slot population"
$ws.Range("N1").Value = "PSU P/Ns"            # was "PSU part number"

# --- column width geometry fix (narrower now that labels are shorter) --
$ws.Columns.Item(6).ColumnWidth = 14.8333333333333    # col F: 18.7109375 -> 15.7109375
$ws.Columns.Item(8).ColumnWidth = 10.8333333333333    # col H: 25.7109375 -> 11.7109375
$ws.Columns.Item(13).ColumnWidth = 12.8333333333333   # col M: 19.7109375 -> 13.7109375
$ws.Columns.Item(14).ColumnWidth = 7.83333333333333   # col N: 15.7109375 -> 8.7109375
